# Collapse duplicate course entries so only one unique "Thwomp Ruins" /
# "Tour Paris Promenade" remain, recomputing meanf (B) and rmkd (D) for the
# affected tracks. Row 43/44 also swap which track name they display,
# because after recomputation "Thwomp Ruins" now sorts ahead of
# "Tour Paris Promenade" in ascending mean.rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 now shows "Thwomp Ruins" with its recomputed stats.
$ws.Range("A43").Value = "Thwomp Ruins"
$ws.Range("B43").Value = 0.72759375775852
$ws.Range("D43").Value = 363.79687887926

# Row 44 now shows "Tour Paris Promenade" with its recomputed stats.
$ws.Range("A44").Value = "Tour Paris Promenade"
$ws.Range("B44").Value = 0.734535215055572
$ws.Range("D44").Value = 367.267607527786

# Other tracks whose meanf / rmkd shifted slightly after the collapse.
$ws.Range("B4").Value = 0.537462078211016
$ws.Range("D4").Value = 268.731039105508

$ws.Range("B11").Value = 0.570796331225673
$ws.Range("D11").Value = 285.398165612837

$ws.Range("B17").Value = 0.5973397048452
$ws.Range("D17").Value = 298.6698524226

$ws.Range("B25").Value = 0.685250221723626
$ws.Range("D25").Value = 342.625110861813

$ws.Range("B29").Value = 0.693760190810109
$ws.Range("D29").Value = 346.880095405054

$ws.Range("B80").Value = 0.974679049688625
$ws.Range("D80").Value = 487.339524844312
